# Apply the day-over-day cryptos price/volume refresh (coinranking.com scrape).
# Column D ("Price") holds numeric-looking strings (dot-grouped thousands, e.g. "34.366.15")
# that must stay text -- a leading apostrophe forces Excel to keep them as literal text
# instead of re-parsing/rounding them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''34.366.15'
$ws.Range("E2").Value = '  +0.82%  '
# Row 3
$ws.Range("D3").Value = '''1.787.82'
$ws.Range("E3").Value = '  +0.41%  '
# Row 4
$ws.Range("E4").Value = '  -0.06%  '
# Row 5
$ws.Range("D5").Value = '''226.13'
$ws.Range("E5").Value = '  +0.31%  '
# Row 6
$ws.Range("D6").Value = '''0.556'
$ws.Range("E6").Value = '  +1.93%  '
# Row 7
$ws.Range("E7").Value = '  -0.06%  '
# Row 8
$ws.Range("D8").Value = '''32.97'
$ws.Range("E8").Value = '  +3.72%  '
# Row 9
$ws.Range("E9").Value = '  +1.14%  '
# Row 10
$ws.Range("E10").Value = '  +0.51%  '
# Row 11
$ws.Range("D11").Value = '''0.0945'
$ws.Range("E11").Value = '  -0.20%  '
# Row 12
$ws.Range("D12").Value = '''2.046.36'
$ws.Range("E12").Value = '  +0.44%  '
# Row 13
$ws.Range("D13").Value = '''11.17'
$ws.Range("E13").Value = '  +2.32%  '
# Row 14
$ws.Range("D14").Value = '''1.768.23'
$ws.Range("E14").Value = '  -0.92%  '
# Row 15
$ws.Range("D15").Value = '''0.635'
$ws.Range("E15").Value = '  +2.27%  '
# Row 16
$ws.Range("D16").Value = '''34.347.63'
$ws.Range("E16").Value = '  +0.75%  '
# Row 17
$ws.Range("E17").Value = '  +2.66%  '
# Row 18
$ws.Range("D18").Value = '''68.42'
$ws.Range("E18").Value = '  +1.32%  '
# Row 19
$ws.Range("D19").Value = '''245.29'
$ws.Range("E19").Value = '  -0.07%  '
# Row 20
$ws.Range("D20").Value = '''0.0₃0795'
$ws.Range("E20").Value = '  +0.98%  '
# Row 21
$ws.Range("D21").Value = '''11.24'
$ws.Range("E21").Value = '  +3.49%  '
# Row 22
$ws.Range("E22").Value = '  -0.22%  '
# Row 23
$ws.Range("D23").Value = '''4.15'
$ws.Range("E23").Value = '  +1.47%  '
# Row 24
$ws.Range("D24").Value = '''168.64'
$ws.Range("E24").Value = '  +4.26%  '
# Row 25
$ws.Range("D25").Value = '''2.06'
$ws.Range("E25").Value = '  +1.93%  '
# Row 26
$ws.Range("E26").Value = '  +3.22%  '
# Row 27
$ws.Range("D27").Value = '''16.56'
$ws.Range("E27").Value = '  +1.97%  '
# Row 28
$ws.Range("E28").Value = '  +1.60%  '
# Row 29
$ws.Range("E29").Value = '  -0.18%  '
# Row 30
$ws.Range("E30").Value = '  +8.03%  '
# Row 31
$ws.Range("D31").Value = '''0.0526'
$ws.Range("E31").Value = '  +1.84%  '
# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.80'
$ws.Range("E32").Value = '  +2.65%  '
# Row 33
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.23'
$ws.Range("E33").Value = '  +0.59%  '
# Row 34
$ws.Range("D34").Value = '''1.82'
$ws.Range("E34").Value = '  +1.52%  '
# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '''2.60'
$ws.Range("E35").Value = '  +6.15%  '
# Row 36
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '''1.409.75'
$ws.Range("E36").Value = '  -2.53%  '
# Row 37
$ws.Range("D37").Value = '''0.682'
$ws.Range("E37").Value = '  +4.92%  '
# Row 38
$ws.Range("E38").Value = '  +2.85%  '
# Row 39
$ws.Range("D39").Value = '''0.0191'
$ws.Range("E39").Value = '  +0.40%  '
# Row 40
$ws.Range("D40").Value = '''84.46'
$ws.Range("E40").Value = '  +4.86%  '
# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.78'
$ws.Range("E41").Value = '  +2.63%  '
# Row 42
$ws.Range("D42").Value = '''2.40'
$ws.Range("E42").Value = '  +0.11%  '
# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '''14.12'
$ws.Range("E43").Value = '  +2.56%  '
# Row 44
$ws.Range("E44").Value = '  +2.60%  '
# Row 45
$ws.Range("E45").Value = '  +2.17%  '
# Row 46
$ws.Range("D46").Value = '''1.10'
$ws.Range("E46").Value = '  +2.59%  '
# Row 47
$ws.Range("D47").Value = '''6.08'
$ws.Range("E47").Value = '  +0.48%  '
# Row 48
$ws.Range("D48").Value = '''1.946.69'
# Row 49
$ws.Range("D49").Value = '''105.42'
$ws.Range("E49").Value = '  +0.99%  '
# Row 50
$ws.Range("E50").Value = '  -0.11%  '
# Row 51
$ws.Range("E51").Value = '  -3.38%  '
